$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 609.6
$ws.Range("I6").Value = 609.6
$ws.Range("K6").Value = 1828.8
$ws.Range("M6").Value = -1716.8
$ws.Range("H11").Value = 408
$ws.Range("I11").Value = 408
$ws.Range("K11").Value = 408
$ws.Range("M11").Value = -268
$ws.Range("H28").Value = 6889.933
$ws.Range("I28").Value = 7346.357
$ws.Range("K28").Value = 7346.357
$ws.Range("M28").Value = -6861.357
$ws.Range("H53").Value = 75
$ws.Range("I53").Value = 75
$ws.Range("K53").Value = 75
$ws.Range("M53").Value = 562
$ws.Range("H100").Value = 5580.7144
$ws.Range("I100").Value = 2021.8334
$ws.Range("K100").Value = 2021.8334
$ws.Range("M100").Value = -1480.8334
$ws.Range("H118").Value = 931.3570999999999
$ws.Range("J118").Value = 450
$ws.Range("L118").Value = 1350
$ws.Range("N118").Value = -4664
$ws.Range("H132").Value = 2851.05
$ws.Range("I132").Value = 2251.6428
$ws.Range("J132").Value = 4249.6665
$ws.Range("K132").Value = 6754.928400000001
$ws.Range("L132").Value = 12748.9995
$ws.Range("M132").Value = -4224.928400000001
$ws.Range("N132").Value = -17808.9995

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1192
$ws.Range("I45").Value = 1192
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1192
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -815
$ws.Range("N45").ClearContents()
$ws.Range("H122").Value = 12348001
$ws.Range("I122").Value = 15875188
$ws.Range("J122").Value = 2848.5
$ws.Range("K122").Value = 47625564
$ws.Range("L122").Value = 8545.5
$ws.Range("M122").Value = -47623114
$ws.Range("N122").Value = -13445.5
$ws.Range("H135").Value = 68732.125
$ws.Range("J135").Value = 68732.125
$ws.Range("L135").Value = 68732.125
$ws.Range("N135").Value = -78872.125

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 935.8182
$ws.Range("I20").Value = 1015.25
$ws.Range("J20").Value = 724
$ws.Range("K20").Value = 1015.25
$ws.Range("L20").Value = 724
$ws.Range("M20").Value = -768.25
$ws.Range("N20").Value = -1218
$ws.Range("I86").Value = 3052.85
$ws.Range("J86").Value = 64876896
$ws.Range("K86").Value = 3052.85
$ws.Range("L86").Value = 64876896
$ws.Range("M86").Value = -1929.85
$ws.Range("N86").Value = -64879142
$ws.Range("I89").Value = 3052.85
$ws.Range("J89").Value = 64876896
$ws.Range("K89").Value = 15264.25
$ws.Range("L89").Value = 324384480
$ws.Range("M89").Value = -9648.25
$ws.Range("N89").Value = -324395712
$ws.Range("H105").Value = 2538.9312
$ws.Range("I105").Value = 2031.5
$ws.Range("K105").Value = 2031.5
$ws.Range("M105").Value = -284.5
$ws.Range("H107").Value = 2163.8975
$ws.Range("I107").Value = 814.55554
$ws.Range("K107").Value = 814.55554
$ws.Range("M107").Value = 1105.44446
$ws.Range("H140").Value = 72500
$ws.Range("J140").Value = 72500
$ws.Range("L140").Value = 72500
$ws.Range("N140").Value = -82860

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3774
$ws.Range("I16").Value = 3612.2144
$ws.Range("K16").Value = 3612.2144
$ws.Range("M16").Value = -3325.2144
$ws.Range("H22").Value = 1222.5294
$ws.Range("I22").Value = 1185.3
$ws.Range("J22").Value = 1275.7142
$ws.Range("K22").Value = 1185.3
$ws.Range("L22").Value = 1275.7142
$ws.Range("M22").Value = -835.3
$ws.Range("N22").Value = -1975.7142
$ws.Range("H31").Value = 3018.611
$ws.Range("I31").Value = 1153.1666
$ws.Range("J31").Value = 6749.5
$ws.Range("K31").Value = 1153.1666
$ws.Range("L31").Value = 6749.5
$ws.Range("M31").Value = -858.1666
$ws.Range("N31").Value = -7339.5
$ws.Range("H34").Value = 3018.611
$ws.Range("I34").Value = 1153.1666
$ws.Range("J34").Value = 6749.5
$ws.Range("K34").Value = 1153.1666
$ws.Range("L34").Value = 6749.5
$ws.Range("M34").Value = -951.1666
$ws.Range("N34").Value = -7153.5
$ws.Range("H64").Value = 75000
$ws.Range("J64").Value = 75000
$ws.Range("L64").Value = 75000
$ws.Range("N64").Value = -75496
$ws.Range("H67").Value = 75000
$ws.Range("J67").Value = 75000
$ws.Range("L67").Value = 75000
$ws.Range("N67").Value = -76716
$ws.Range("H113").Value = 3774
$ws.Range("I113").Value = 3612.2144
$ws.Range("K113").Value = 3612.2144
$ws.Range("M113").Value = -1442.2144
$ws.Range("H132").Value = 8599.6
$ws.Range("I132").Value = 8599.6
$ws.Range("K132").Value = 25798.8
$ws.Range("M132").Value = -23268.8

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 97.17391000000001
$ws.Range("I2").Value = 118.69231
$ws.Range("J2").Value = 69.2
$ws.Range("K2").Value = 712.15386
$ws.Range("L2").Value = 415.2
$ws.Range("M2").Value = -599.15386
$ws.Range("N2").Value = -641.2
$ws.Range("H38").Value = 866.2308
$ws.Range("I38").Value = 30
$ws.Range("K38").Value = 90
$ws.Range("M38").Value = 257
$ws.Range("H44").Value = 1335.4546
$ws.Range("J44").Value = 1288.75
$ws.Range("L44").Value = 3866.25
$ws.Range("N44").Value = -4662.25
$ws.Range("H139").Value = 5751.8125
$ws.Range("I139").Value = 7338.1665
$ws.Range("J139").Value = 4800
$ws.Range("K139").Value = 22014.4995
$ws.Range("L139").Value = 14400
$ws.Range("M139").Value = -16874.4995
$ws.Range("N139").Value = -24680

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 30001
$ws.Range("J52").Value = 30001
$ws.Range("L52").Value = 30001
$ws.Range("N52").Value = -30519
$ws.Range("H54").Value = 7000
$ws.Range("J54").Value = 7000
$ws.Range("L54").Value = 7000
$ws.Range("N54").Value = -7780
$ws.Range("H97").Value = 9999.727999999999
$ws.Range("I97").Value = 814.4286
$ws.Range("K97").Value = 814.4286
$ws.Range("M97").Value = -318.4286
$ws.Range("H113").Value = 30308662
$ws.Range("I113").Value = 52634148
$ws.Range("J113").Value = 9785.714
$ws.Range("K113").Value = 52634148
$ws.Range("L113").Value = 9785.714
$ws.Range("M113").Value = -52631978
$ws.Range("N113").Value = -14125.714

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3698.45
$ws.Range("I22").Value = 3568.7144
$ws.Range("J22").Value = 3768.3076
$ws.Range("K22").Value = 3568.7144
$ws.Range("L22").Value = 3768.3076
$ws.Range("M22").Value = -3273.7144
$ws.Range("N22").Value = -4358.3076
$ws.Range("H27").Value = 3698.45
$ws.Range("I27").Value = 3568.7144
$ws.Range("J27").Value = 3768.3076
$ws.Range("K27").Value = 3568.7144
$ws.Range("L27").Value = 3768.3076
$ws.Range("M27").Value = -3461.7144
$ws.Range("N27").Value = -3982.3076
$ws.Range("H46").Value = 873.90625
$ws.Range("I46").Value = 672.9231
$ws.Range("J46").Value = 1744.8334
$ws.Range("K46").Value = 672.9231
$ws.Range("L46").Value = 1744.8334
$ws.Range("M46").Value = -484.9231
$ws.Range("N46").Value = -2120.8334
$ws.Range("H55").Value = 1003.7
$ws.Range("I55").Value = 174.33333
$ws.Range("J55").Value = 2247.75
$ws.Range("K55").Value = 174.33333
$ws.Range("L55").Value = 2247.75
$ws.Range("M55").Value = -1.333329999999989
$ws.Range("N55").Value = -2593.75
$ws.Range("H61").Value = 3001.6453
$ws.Range("I61").Value = 1801.25
$ws.Range("K61").Value = 1801.25
$ws.Range("M61").Value = -1599.25
$ws.Range("H113").Value = 3001.6453
$ws.Range("I113").Value = 1801.25
$ws.Range("K113").Value = 1801.25
$ws.Range("M113").Value = 368.75

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 92249.5
$ws.Range("J16").Value = 92249.5
$ws.Range("L16").Value = 92249.5
$ws.Range("N16").Value = -92833.5
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H121").Value = 74449.5
$ws.Range("J121").Value = 74449.5
$ws.Range("L121").Value = 74449.5
$ws.Range("N121").Value = -77943.5
